$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# --- Update "series" sheet (G column) costs ------------------------------
$ws = $wb.Worksheets.Item("series")
$newG = @(0.16,0.12,0.19,0.19,0.18,0.18,0.1,0.19,0.19,0.2,0.13,0.18,0.16,0.12,0.1,0.17,0.14,0.14,0.13,0.19,0.13,0.2,0.19,0.19)
for ($i = 0; $i -lt $newG.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $newG[$i]
}

# --- Remove the obsolete "elements" sheet ---------------------------------
$elements = $wb.Worksheets.Item("elements")
$elements.Delete() | Out-Null

# --- Restore per-sheet cursor positions & active tab ----------------------
$wsSeries = $wb.Worksheets.Item("series")
$wsSeries.Range("G4").Select() | Out-Null

$wsElementsTest = $wb.Worksheets.Item("elements test")
$wsElementsTest.Range("I2").Select() | Out-Null

$wsConectElectric = $wb.Worksheets.Item("conect_electric")
$wsConectElectric.Range("D5").Select() | Out-Null

$wsConectThermal = $wb.Worksheets.Item("conect_thermal")
$wsConectThermal.Range("D1").Select() | Out-Null

# "elements test" ends up as the active/selected tab
$wsElementsTest.Activate()
$wsElementsTest.Range("I2").Select() | Out-Null
